$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- New rows of content on Feuil1, continuing the "methodes" table ---
# Shared strings must be created in the same order they first appear in the
# target workbook, so the cell values below are written in that order:
#   B60, A60, A61, B61, A63, A64, B63, B64, A66, B66, A67, B67

$ws1.Range("B60").Value = "inner join avec MEMBRE pour recuperer l'image du profile du MEMBRE qui a envoye  le message"
$ws1.Range("A60").Value = "AJOUTER un attribut String IMAGEPATH dans le bean MESSAGE"

$ws1.Range("A61").Value = "AJOUTER un attribut String PSEUDO dans le bean MESSAGE"
$ws1.Range("B61").Value = "inner join avec MEMBRE pour recuperer le Pseudo du profile du MEMBRE qui a envoye  le message"

$ws1.Range("A63").Value = "AJOUTER un attribut String IMAGEPATH dans le bean CLINSDOEIL"
$ws1.Range("A64").Value = "AJOUTER un attribut String PSEUDO dans le bean  CLINSDOEIL"
$ws1.Range("B63").Value = "inner join avec MEMBRE pour recuperer l'image du profile du MEMBRE qui a envoye  le  CLINSDOEIL"
$ws1.Range("B64").Value = "inner join avec MEMBRE pour recuperer le Pseudo du profile du MEMBRE qui a envoye  le CLINSDOEIL"

$ws1.Range("A66").Value = "get toute la liste de niveux de membre"
$ws1.Range("B66").Value = "return= objet niveaumembre"

$ws1.Range("A67").Value = "get membre by id"
$ws1.Range("B67").Value = "MembreManager.getMembreById(int memberid)"

# --- Formatting: column A (and some column B) cells use the existing red
# "note" font style already used throughout this table (same as A58/A60 etc.) ---
$ws1.Range("A60").Font.Color = 255
$ws1.Range("A61").Font.Color = 255
$ws1.Range("A63").Font.Color = 255
$ws1.Range("A64").Font.Color = 255

$ws1.Range("A66").Font.Color = 255
$ws1.Range("B66").Font.Color = 255
$ws1.Range("A67").Font.Color = 255
$ws1.Range("B67").Font.Color = 255

# Row 68 stays empty but keeps the same red styling carried over.
$ws1.Range("A68").Font.Color = 255
$ws1.Range("B68").Font.Color = 255

# --- Selection / view state: the sheet was left scrolled down with A67:B67
# selected on Feuil1 (which remains the active sheet/tab). ---
$ws1.Range("A67:B67").Select() | Out-Null
